# Insert a new data row at row 23 (pushing existing rows 23-112 down to 24-113)
# and populate it with a new "Choclero / Primera" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).EntireRow.Insert()

$ws.Range("A23").Value = 2
$ws.Range("B23").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C23").Value = "Coquimbo"
$ws.Range("D23").Value = 44608
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = 100112024
$ws.Range("G23").Value = "Choclo"
$ws.Range("H23").Value = "Choclero"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 70000
$ws.Range("K23").Value = 130
$ws.Range("L23").Value = 150
$ws.Range("M23").Value = 140
$ws.Range("N23").Value = "$/unidad"
$ws.Range("O23").Value = "Provincia de Limarí"
$ws.Range("P23").Value = 140
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = "Hortaliza"
